$wb = $excel.ActiveWorkbook

# --- 毒圈 (poison circle) sheet: update the D column (移动/move) values for rows 4-11.
# All other changed cells on this sheet (G,H,I,J) are formulas that recompute
# automatically once D changes.
$ws5 = $wb.Worksheets.Item("毒圈")

$ws5.Cells.Item(4, 4).Value = 180
$ws5.Cells.Item(5, 4).Value = 160
$ws5.Cells.Item(6, 4).Value = 140
$ws5.Cells.Item(7, 4).Value = 120
$ws5.Cells.Item(8, 4).Value = 100
$ws5.Cells.Item(9, 4).Value = 80
$ws5.Cells.Item(10, 4).Value = 60
$ws5.Cells.Item(11, 4).Value = 40

# --- Switch the active sheet/tab from 角色 (character) to 毒圈 (poison circle),
# matching the author moving on to debug the circle data.
$ws4 = $wb.Worksheets.Item("角色")
$ws4.Select()
$ws5.Select()
$ws5.Range("G17").Select()
